# Stash the real username/password in place of the placeholder tokens that
# were being used for the login automation steps, and leave the active
# selection on the next row (C12) where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "gracz"
$ws.Range("C11").Value = "passwerd"

$ws.Range("C12").Select()
